$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-25 Friday", "2025-07-26 Saturday"),
    @("170÷4=", "923÷2="),
    @("381÷3=", "183÷4="),
    @("602÷2=", "127÷7="),
    @("128÷3=", "250÷9="),
    @("339÷8=", "442÷3="),
    @("559÷5=", "757÷4="),
    @("916÷2=", "303÷2="),
    @("268÷7=", "310÷6="),
    @("718÷2=", "878÷9="),
    @("305÷3=", "834÷2="),
    @("320÷8=", "244÷2="),
    @("115÷4=", "642÷2="),
    @("314÷8=", "895÷6="),
    @("715÷7=", "159÷7="),
    @("771÷6=", "482÷2="),
    @("563÷4=", "709÷2="),
    @("338÷4=", "702÷6="),
    @("924÷4=", "280÷3="),
    @("343÷4=", "852÷5="),
    @("860÷6=", "264÷2="),
    @("362÷8=", "310÷8="),
    @("189÷2=", "477÷6="),
    @("308÷2=", "479÷8="),
    @("956÷6=", "731÷7="),
    @("520÷8=", "519÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
